# Adds a friendly-text second argument to the HYPERLINK() formulas in S2:Y4.
# Note: per the source diff, only column S gets a well-formed two-argument
# HYPERLINK(url; text) call. Columns T:Y are updated with the same literal
# (buggy) replacement pattern as the original commit, which drops the
# closing quote after the URL before appending `; "label")` - i.e. the
# resulting formula text is reproduced verbatim, quirks included.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/artfynd/A 30840-2023.xlsx"; "A 30840-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/kartor/A 30840-2023.png; "A 30840-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/knärot/A 30840-2023.png; "A 30840-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomål/A 30840-2023.docx; "A 30840-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomålsmail/A 30840-2023.docx; "A 30840-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsyn/A 30840-2023.docx; "A 30840-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsynsmail/A 30840-2023.docx; "A 30840-2023")'
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/artfynd/A 30841-2023.xlsx"; "A 30841-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/kartor/A 30841-2023.png; "A 30841-2023")'
$ws.Range("U3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/knärot/A 30841-2023.png; "A 30841-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomål/A 30841-2023.docx; "A 30841-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomålsmail/A 30841-2023.docx; "A 30841-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsyn/A 30841-2023.docx; "A 30841-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsynsmail/A 30841-2023.docx; "A 30841-2023")'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/artfynd/A 30839-2023.xlsx"; "A 30839-2023")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/kartor/A 30839-2023.png; "A 30839-2023")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomål/A 30839-2023.docx; "A 30839-2023")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomålsmail/A 30839-2023.docx; "A 30839-2023")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsyn/A 30839-2023.docx; "A 30839-2023")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsynsmail/A 30839-2023.docx; "A 30839-2023")'
